$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.999.51"
$ws.Range("E2").Value = "'  -1.98%  "
$ws.Range("D3").Value = "'3.122.31"
$ws.Range("E3").Value = "'  -0.60%  "
$ws.Range("E4").Value = "'  +0.10%  "
$ws.Range("D5").Value = "'596.42"
$ws.Range("E5").Value = "'  -2.22%  "
$ws.Range("D6").Value = "'136.08"
$ws.Range("E6").Value = "'  -5.24%  "
$ws.Range("E7").Value = "'  +0.02%  "
$ws.Range("D8").Value = "'3.114.93"
$ws.Range("E8").Value = "'  -0.72%  "
$ws.Range("D9").Value = "'0.517"
$ws.Range("E9").Value = "'  -2.47%  "
$ws.Range("D10").Value = "'0.145"
$ws.Range("E10").Value = "'  -3.70%  "
$ws.Range("D11").Value = "'5.19"
$ws.Range("E11").Value = "'  -4.34%  "
$ws.Range("D12").Value = "'0.459"
$ws.Range("E12").Value = "'  -3.63%  "
$ws.Range("D13").Value = "'0.0000247"
$ws.Range("E13").Value = "'  -3.63%  "
$ws.Range("D14").Value = "'34.25"
$ws.Range("E14").Value = "'  -3.84%  "
$ws.Range("D15").Value = "'3.630.26"
$ws.Range("E15").Value = "'  -0.57%  "
$ws.Range("E16").Value = "'  +1.51%  "
$ws.Range("D17").Value = "'62.960.89"
$ws.Range("E17").Value = "'  -1.97%  "
$ws.Range("D18").Value = "'3.123.34"
$ws.Range("E18").Value = "'  -0.52%  "
$ws.Range("D19").Value = "'6.73"
$ws.Range("E19").Value = "'  -2.23%  "
$ws.Range("D20").Value = "'476.60"
$ws.Range("E20").Value = "'  -0.22%  "
$ws.Range("D21").Value = "'14.16"
$ws.Range("E21").Value = "'  -3.88%  "
$ws.Range("D22").Value = "'0.697"
$ws.Range("E22").Value = "'  -3.86%  "
$ws.Range("D23").Value = "'7.66"
$ws.Range("E23").Value = "'  -2.34%  "
$ws.Range("D24").Value = "'87.45"
$ws.Range("E24").Value = "'  +2.34%  "
$ws.Range("D25").Value = "'13.00"
$ws.Range("E25").Value = "'  -5.07%  "
$ws.Range("E26").Value = "'  +0.18%  "
$ws.Range("D27").Value = "'2.71"
$ws.Range("E27").Value = "'  -2.57%  "
$ws.Range("D28").Value = "'7.22"
$ws.Range("E28").Value = "'  -2.84%  "
$ws.Range("D29").Value = "'7.91"
$ws.Range("E29").Value = "'  -7.93%  "
$ws.Range("D30").Value = "'2.08"
$ws.Range("E30").Value = "'  -0.36%  "
$ws.Range("D31").Value = "'27.09"
$ws.Range("E31").Value = "'  +1.73%  "
$ws.Range("E32").Value = "'  -0.04%  "
$ws.Range("D33").Value = "'0.107"
$ws.Range("E33").Value = "'  -8.55%  "
$ws.Range("D34").Value = "'2.53"
$ws.Range("E34").Value = "'  -4.29%  "
$ws.Range("E35").Value = "'  -2.75%  "
$ws.Range("D36").Value = "'5.82"
$ws.Range("E36").Value = "'  -2.43%  "
$ws.Range("D37").Value = "'51.94"
$ws.Range("E37").Value = "'  -0.87%  "
$ws.Range("D38").Value = "'0.0₃0711"
$ws.Range("E38").Value = "'  -4.91%  "
$ws.Range("D39").Value = "'0.0387"
$ws.Range("E39").Value = "'  -2.40%  "
$ws.Range("D40").Value = "'421.66"
$ws.Range("E40").Value = "'  -7.22%  "
$ws.Range("E41").Value = "'  -0.68%  "
$ws.Range("D42").Value = "'8.28"
$ws.Range("E42").Value = "'  -0.86%  "
$ws.Range("D43").Value = "'2.66"
$ws.Range("E43").Value = "'  -11.95%  "
$ws.Range("D44").Value = "'2.872.99"
$ws.Range("E44").Value = "'  -0.19%  "
$ws.Range("D45").Value = "'0.264"
$ws.Range("E45").Value = "'  -0.43%  "
$ws.Range("B46").Value = "'USDe"
$ws.Range("C46").Value = "'https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D46").Value = "'0.999"
$ws.Range("E46").Value = "'  -0.02%  "
$ws.Range("B47").Value = "'Fetch.AI"
$ws.Range("C47").Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").Value = "'2.13"
$ws.Range("E47").Value = "'  -6.13%  "
$ws.Range("D48").Value = "'25.73"
$ws.Range("E48").Value = "'  -3.14%  "
$ws.Range("D49").Value = "'0.113"
$ws.Range("E49").Value = "'  -0.91%  "
$ws.Range("D50").Value = "'2.27"
$ws.Range("E50").Value = "'  -7.09%  "
$ws.Range("D51").Value = "'118.42"
$ws.Range("E51").Value = "'  -2.08%  "
